# Update Name of Algo
# Apply updated KNN imputation results to specific cells in column A and D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = -7.444
$ws.Range("A12").Value = -21.626
$ws.Range("D14").Value = -7.391
$ws.Range("D26").Value = -8.000999999999999
$ws.Range("A27").Value = -21.718
$ws.Range("D31").Value = -8.218999999999999
$ws.Range("A32").Value = -21.595
$ws.Range("D35").Value = -7.939
$ws.Range("A36").Value = -20.339
$ws.Range("D37").Value = -7.741
$ws.Range("A38").Value = -19.741
$ws.Range("D45").Value = -7.539
$ws.Range("A46").Value = -21.789
$ws.Range("D52").Value = -7.87
$ws.Range("A54").Value = -21.703
$ws.Range("A55").Value = -22.196
$ws.Range("A56").Value = -22.095
$ws.Range("D57").Value = -8.289999999999999
$ws.Range("A67").Value = -21.565
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.445
$ws.Range("D81").Value = -7.309
$ws.Range("A83").Value = -21.776
$ws.Range("D83").Value = -8.289
$ws.Range("A86").Value = -22.035
$ws.Range("A91").Value = -21.564
$ws.Range("A93").Value = -21.259
$ws.Range("A99").Value = -20.43
$ws.Range("D100").Value = -8.238
$ws.Range("D102").Value = -7.752
